# chore: adapt column header formatting to respective input file names
#
# - Rename the "_old" column headers (A1:J1) to use the "_FV2310" suffix.
# - Rename the "_new" column headers (L1:U1) to use the "_FV2404" suffix.
# - Leave the "diff" column header (K1) untouched.
# - Turn the A1:U57 range into a proper Excel Table ("Table1") with
#   autofilter on the header row.
# - Freeze the header row (split/freeze after row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffix = "_old"
$newSuffix = "_new"
$fv2310Suffix = "_FV2310"
$fv2404Suffix = "_FV2404"

$lastDataRow = 57
$firstOldCol = 1    # A
$lastOldCol = 10    # J
$firstNewCol = 12   # L
$lastNewCol = 21    # U

# Columns A..J: "<name>_old" -> "<name>_FV2310"
for ($col = $firstOldCol; $col -le $lastOldCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = ($cell.Value2).Replace($oldSuffix, $fv2310Suffix)
}

# Column K ("diff") is intentionally left as-is.

# Columns L..U: "<name>_new" -> "<name>_FV2404"
for ($col = $firstNewCol; $col -le $lastNewCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = ($cell.Value2).Replace($newSuffix, $fv2404Suffix)
}

# Turn the full data range (including the freshly renamed header row) into
# an Excel Table so the column headers and autofilter are wired up.
$tableRange = $ws.Range("A1:U" + $lastDataRow)
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"

# Freeze the header row: select A2 then freeze so row 1 stays visible.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
